# Update LR-pairs data with new TPM-derived values (16 Sending x Target cluster combinations)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$rows = @(
  @("ECs", "ECs", 3, 1, 6.449754000000001, 19.349262, 0.03479900749229446, 0.03479900749229446, 3, 1, 237.3377026666667, 712.013108, 0.9887685707142667, 0.9887685707142667, 1530.769797125144, 13776.9281741263, 0.03440816490043105, 0.03440816490043105),
  @("ECs", "FAPs", 3, 1, 6.449754000000001, 19.349262, 0.03479900749229446, 0.03479900749229446, 3, 1, 1.025352333333333, 3.076057, 0.004271702935173513, 0.004271702935173513, 6.613270313326002, 59.51943281993401, 0.0001486510224459593, 0.0001486510224459593),
  @("ECs", "MuSCs", 3, 1, 6.449754000000001, 19.349262, 0.03479900749229446, 0.03479900749229446, 3, 1, 1.587950666666667, 4.763852, 0.006615534293133127, 0.006615534293133127, 10.241891164136, 92.17702047722402, 0.0002302140274322706, 0.0002302140274322706),
  @("ECs", "Resolving-Mac", 3, 1, 6.449754000000001, 19.349262, 0.03479900749229446, 0.03479900749229446, 3, 1, 0.08261766666666666, 0.247853, 0.0003441920574266213, 0.0003441920574266213, 0.5328636260540001, 4.795772634486001, 0.00001197754198517724, 0.00001197754198517724),
  @("FAPs", "ECs", 3, 1, 123.027733, 369.083199, 0.663783921437469, 0.6637839214374691, 3, 1, 237.3377026666667, 712.013108, 0.9887685707142667, 0.9887685707142667, 29199.11951450806, 262792.0756305724, 0.6563286792628373, 0.6563286792628374),
  @("FAPs", "FAPs", 3, 1, 123.027733, 369.083199, 0.663783921437469, 0.6637839214374691, 3, 1, 1.025352333333333, 3.076057, 0.004271702935173513, 0.004271702935173513, 126.1467730962603, 1135.320957866343, 0.002835487725525421, 0.002835487725525421),
  @("FAPs", "MuSCs", 3, 1, 123.027733, 369.083199, 0.663783921437469, 0.6637839214374691, 3, 1, 1.587950666666667, 4.763852, 0.006615534293133127, 0.006615534293133127, 195.3619706358387, 1758.257735722548, 0.004391285295499962, 0.004391285295499962),
  @("FAPs", "Resolving-Mac", 3, 1, 123.027733, 369.083199, 0.663783921437469, 0.6637839214374691, 3, 1, 0.08261766666666666, 0.247853, 0.0003441920574266213, 0.0003441920574266213, 10.16426423574967, 91.47837812174699, 0.0002284691536062732, 0.0002284691536062733),
  @("MuSCs", "ECs", 3, 1, 55.79038633333334, 167.371159, 0.3010114916028843, 0.3010114916028843, 3, 1, 237.3377026666667, 712.013108, 0.9887685707142667, 0.9887685707142667, 13241.16212323913, 119170.4591091522, 0.2976307023207534, 0.2976307023207534),
  @("MuSCs", "FAPs", 3, 1, 55.79038633333334, 167.371159, 0.3010114916028843, 0.3010114916028843, 3, 1, 1.025352333333333, 3.076057, 0.004271702935173513, 0.004271702935173513, 57.20480280445145, 514.843225240063, 0.001285831672200998, 0.001285831672200998),
  @("MuSCs", "MuSCs", 3, 1, 55.79038633333334, 167.371159, 0.3010114916028843, 0.3010114916028843, 3, 1, 1.587950666666667, 4.763852, 0.006615534293133127, 0.006615534293133127, 88.59238117160757, 797.331430544468, 0.001991351845326035, 0.001991351845326035),
  @("MuSCs", "Resolving-Mac", 3, 1, 55.79038633333334, 167.371159, 0.3010114916028843, 0.3010114916028843, 3, 1, 0.08261766666666666, 0.247853, 0.0003441920574266213, 0.0003441920574266213, 4.609271541291889, 41.483443871627, 0.0001036057646038529, 0.0001036057646038529),
  @("Resolving-Mac", "ECs", 2, 0.6666666666666666, 0.07517133333333333, 0.225514, 0.0004055794673521549, 0.000405579467352155, 3, 1, 237.3377026666667, 712.013108, 0.9887685707142667, 0.9887685707142667, 17.84099155972356, 160.568924037512, 0.0004010242302448438, 0.0004010242302448438),
  @("Resolving-Mac", "FAPs", 2, 0.6666666666666666, 0.07517133333333333, 0.225514, 0.0004055794673521549, 0.000405579467352155, 3, 1, 1.025352333333333, 3.076057, 0.004271702935173513, 0.004271702935173513, 0.0770771020331111, 0.693693918298, 0.00000173251500113431, 0.00000173251500113431),
  @("Resolving-Mac", "MuSCs", 2, 0.6666666666666666, 0.07517133333333333, 0.225514, 0.0004055794673521549, 0.000405579467352155, 3, 1, 1.587950666666667, 4.763852, 0.006615534293133127, 0.006615534293133127, 0.1193683688808889, 1.074315319928, 0.000002683124874858848, 0.000002683124874858849),
  @("Resolving-Mac", "Resolving-Mac", 2, 0.6666666666666666, 0.07517133333333333, 0.225514, 0.0004055794673521549, 0.000405579467352155, 3, 1, 0.08261766666666666, 0.247853, 0.0003441920574266213, 0.0003441920574266213, 0.006210480160222221, 0.055894321442, 0.0000001395972313179314, 0.0000001395972313179314)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
  $r = $startRow + $i
  $row = $rows[$i]
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = "Tnc"
  $ws.Cells.Item($r, 3).Value = "Ptprb"
  $ws.Cells.Item($r, 4).Value = $row[1]
  for ($c = 2; $c -lt $row.Count; $c++) {
    $ws.Cells.Item($r, 3 + $c).Value = $row[$c]
  }
}

